$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The paragraph that used to live in row 6 ("We may be the United Nations...
# ...Climate Change and International Terrorism...") is being split into two
# separate rows. Insert a fresh row at position 7 to make room for the
# second half; everything from the old row 7 onward shifts down by one.
$ws.Rows(7).Insert()

# Row 6: first half of the split paragraph (unchanged "no" answer).
$ws.Range("A6").Value = "`nWe may be the United Nations of Great Grief but our love, strength in diversity, tolerance and unity is `nfar more powerful than their bombs and guns can ever be! `n "
$ws.Range("B6").Value = "no"

# Row 7 (newly inserted): second half of the split paragraph, now coded as
# a relevant ("yes") passage with its own classification.
$ws.Range("A7").Value = "`nClimate Change and Intern ational Terrorism are two extremely serious issues which threatens Humanity `nand which must be dealt with all possible means at our collective disposal."
$ws.Range("B7").Value = "yes"
$ws.Range("C7").Value = "action"
$ws.Range("D7").Value = "n.a."
$ws.Range("E7").Value = "global"
$ws.Range("F7").Value = "n.a."
$ws.Range("G7").Value = "utilitarian"
$ws.Range("H7").Value = "Urge to take on action to counter threat to humanity, thus in benefit of all. "

# Row 11 (formerly row 10, "It is therefore...") - the Principle and 30-word
# explanation were revised.
$ws.Range("G11").Value = "egalitarain, utilitarian"
$ws.Range("H11").Value = "Egalitarian in urging all countries, for the reason to save the planet, both utilitarain and egalitarain "

# Row 13 (formerly row 12, "Let us seize this historic opportunity...") - the
# 30-word explanation was revised.
$ws.Range("H13").Value = "Urge to take on action for the benefit of all. "

# Row 14 (formerly row 13, "I am sure that you all agree...") - no longer
# coded as relevant; its Topic/Unit/Scale/Time/Principle/explanation are
# cleared and Relevance flips to "no".
$ws.Range("B14").Value = "no"
$ws.Range("C14:H14").ClearContents()

# Leave the selection where the author apparently left off editing.
$ws.Range("C13").Select() | Out-Null
